# Update of metadata obtained on 6 April 2016: the "Nivel estudios" column
# (column D) was mis-tagged as a measure instead of a dimension, and its
# datatype annotation was wrong. Also add a new metadata row pointing to the
# external mapping workbook for this dimension's codes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3: iaest-measure:nivel-estudios -> iaest-dimension:nivel-estudios
$ws.Range("D3").Value = "iaest-dimension:nivel-estudios"

# D4: medida -> dim
$ws.Range("D4").Value = "dim"

# D5: xsd:string -> skos:Concept
$ws.Range("D5").Value = "skos:Concept"

# New row 6: reference to the external mapping file for this dimension.
# Copy D5's formatting first so the new cell matches the rest of the column.
$ws.Range("D5").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D6").Value = "mapping-nivel-estudios.xlsx"
